$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.042.20"
$ws.Range("E2").Value = "  -2.46%  "
$ws.Range("D3").Value = "1.819.83"
$ws.Range("E3").Value = "  -1.52%  "
$ws.Range("D4").Value = "'1.000"
$ws.Range("E4").Value = "  -1.25%  "
$ws.Range("D5").Value = "'311.01"
$ws.Range("D6").Value = "'1.000"
$ws.Range("E6").Value = "  -1.06%  "
$ws.Range("D7").Value = "'0.4217"
$ws.Range("E7").Value = "  -2.24%  "
$ws.Range("E8").Value = "  -2.13%  "
$ws.Range("D9").Value = "'0.07202"
$ws.Range("E9").Value = "  -2.03%  "
$ws.Range("D10").Value = "'0.8391"
$ws.Range("E10").Value = "  -4.34%  "
$ws.Range("E11").Value = "  -4.02%  "
$ws.Range("D12").Value = "1.811.88"
$ws.Range("E12").Value = "  -1.96%  "
$ws.Range("D13").Value = "'6.642"
$ws.Range("E13").Value = "  -1.23%  "
$ws.Range("D14").Value = "'0.07085"
$ws.Range("E14").Value = "  -0.80%  "
$ws.Range("D15").Value = "'5.273"
$ws.Range("E15").Value = "  -3.16%  "
$ws.Range("D16").Value = "'89.37"
$ws.Range("E16").Value = "  +0.28%  "
$ws.Range("E17").Value = "  -1.31%  "
$ws.Range("D18").Value = "'0.000008781"
$ws.Range("E18").Value = "  -2.48%  "
$ws.Range("D19").Value = "'1.000"
$ws.Range("E19").Value = "  -1.07%  "
$ws.Range("E20").Value = "  -3.59%  "
$ws.Range("D21").Value = "27.009.43"
$ws.Range("E21").Value = "  -2.56%  "
$ws.Range("D22").Value = "'5.118"
$ws.Range("E22").Value = "  -1.93%  "
$ws.Range("D23").Value = "'10.82"
$ws.Range("E23").Value = "  -2.33%  "
$ws.Range("D24").Value = "2.035.46"
$ws.Range("E24").Value = "  -2.00%  "
$ws.Range("D25").Value = "'1.971"
$ws.Range("E25").Value = "  -1.62%  "
$ws.Range("D26").Value = "'151.60"
$ws.Range("E26").Value = "  -2.32%  "
$ws.Range("D27").Value = "'2.218"
$ws.Range("E27").Value = "  +2.07%  "
$ws.Range("E28").Value = "  -2.31%  "
$ws.Range("D29").Value = "'5.215"
$ws.Range("E29").Value = "  -3.42%  "
$ws.Range("D30").Value = "'115.98"
$ws.Range("E30").Value = "  -2.72%  "
$ws.Range("D31").Value = "'0.08743"
$ws.Range("E31").Value = "  -2.32%  "
$ws.Range("D32").Value = "'1.178"
$ws.Range("E32").Value = "  -4.53%  "
$ws.Range("D33").Value = "'2.953"
$ws.Range("E33").Value = "  +1.34%  "
$ws.Range("D34").Value = "'0.7373"
$ws.Range("E34").Value = "  -5.19%  "
$ws.Range("D35").Value = "'4.409"
$ws.Range("E35").Value = "  -3.33%  "
$ws.Range("D36").Value = "'0.9999"
$ws.Range("E36").Value = "  -1.22%  "
$ws.Range("D37").Value = "'1.091"
$ws.Range("E37").Value = "  -4.00%  "
$ws.Range("D38").Value = "'0.01948"
$ws.Range("E38").Value = "  -1.47%  "
$ws.Range("D39").Value = "'0.05229"
$ws.Range("E39").Value = "  -2.15%  "
$ws.Range("D40").Value = "'7.288"
$ws.Range("E40").Value = "  -0.29%  "
$ws.Range("D41").Value = "'2.875"
$ws.Range("E41").Value = "  -0.79%  "
$ws.Range("D42").Value = "'0.1685"
$ws.Range("E42").Value = "  -0.18%  "
$ws.Range("D43").Value = "'0.5025"
$ws.Range("E43").Value = "  -2.14%  "
$ws.Range("D44").Value = "'8.577"
$ws.Range("E44").Value = "  -2.76%  "
$ws.Range("D45").Value = "'10.49"
$ws.Range("E45").Value = "  -1.99%  "
$ws.Range("D46").Value = "'106.06"
$ws.Range("E46").Value = "  -2.61%  "
$ws.Range("D47").Value = "'0.4704"
$ws.Range("E47").Value = "  -1.36%  "
$ws.Range("D48").Value = "'0.9997"
$ws.Range("E48").Value = "  -1.21%  "
$ws.Range("D49").Value = "'0.06349"
$ws.Range("E49").Value = "  -1.95%  "
$ws.Range("B50").Value = "RenderToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D50").Value = "'1.886"
$ws.Range("E50").Value = "  +1.64%  "
$ws.Range("B51").Value = "NEARProtocol"
$ws.Range("C51").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D51").Value = "'1.643"
$ws.Range("E51").Value = "  -2.89%  "
